$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G9").Value = 2.4
$ws.Range("I9").Value = 3.5
$ws.Range("Z9").Value = 8.5
$ws.Range("AE9").Value = 11
$ws.Range("L15").Value = 1.22
$ws.Range("M15").Value = 4
$ws.Range("N15").Value = 1.8
$ws.Range("O15").Value = 2
$ws.Range("G16").Value = 1.62
$ws.Range("H16").Value = 3.55
$ws.Range("I16").Value = 5.7
$ws.Range("J16").Value = 1.09
$ws.Range("K16").Value = 6.5
$ws.Range("L16").Value = 1.42
$ws.Range("M16").Value = 2.7
$ws.Range("N16").Value = 2.25
$ws.Range("O16").Value = 1.6
$ws.Range("P16").Value = 1.47
$ws.Range("Q16").Value = 2.55
$ws.Range("R16").Value = 2.22
$ws.Range("S16").Value = 1.6
$ws.Range("T16").Value = 5.2
$ws.Range("U16").Value = 6.8
$ws.Range("V16").Value = 9.25
$ws.Range("W16").Value = 12.5
$ws.Range("X16").Value = 17
$ws.Range("Y16").Value = 45
$ws.Range("Z16").Value = 6.5
$ws.Range("AA16").Value = 7.5
$ws.Range("AC16").Value = 175
$ws.Range("AE16").Value = 11.25
$ws.Range("AF16").Value = 35
$ws.Range("AG16").Value = 21
$ws.Range("AH16").Value = 150
$ws.Range("AI16").Value = 80
$ws.Range("G17").Value = 1.82
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 4.15
$ws.Range("J17").Value = 1.07
$ws.Range("K17").Value = 7.4
$ws.Range("L17").Value = 1.35
$ws.Range("M17").Value = 3.05
$ws.Range("N17").Value = 2.05
$ws.Range("O17").Value = 1.72
$ws.Range("P17").Value = 1.42
$ws.Range("Q17").Value = 2.75
$ws.Range("T17").Value = 6.3
$ws.Range("U17").Value = 8.5
$ws.Range("V17").Value = 9
$ws.Range("W17").Value = 15.5
$ws.Range("X17").Value = 16.5
$ws.Range("Y17").Value = 35
$ws.Range("Z17").Value = 7.4
$ws.Range("AA17").Value = 7.2
$ws.Range("AB17").Value = 18.5
$ws.Range("AE17").Value = 10.25
$ws.Range("AF17").Value = 23
$ws.Range("AG17").Value = 15
$ws.Range("AH17").Value = 70
$ws.Range("AI17").Value = 45
$ws.Range("AJ17").Value = 55
$ws.Range("N20").Value = 1.75
$ws.Range("O20").Value = 2.05
$ws.Range("H44").Value = 3
$ws.Range("U44").Value = 11.5
$ws.Range("V44").Value = 9.75
$ws.Range("AB44").Value = 14.5
$ws.Range("AC44").Value = 80
$ws.Range("AH44").Value = 17
$ws.Range("G53").Value = 1.88
$ws.Range("H53").Value = 3.55
$ws.Range("I53").Value = 3.5
$ws.Range("J53").Value = 1.06
$ws.Range("K53").Value = 7.5
$ws.Range("L53").Value = 1.29
$ws.Range("M53").Value = 3.3
$ws.Range("N53").Value = 1.85
$ws.Range("O53").Value = 1.85
$ws.Range("P53").Value = 1.4
$ws.Range("Q53").Value = 2.75
$ws.Range("R53").Value = 1.78
$ws.Range("S53").Value = 1.93
$ws.Range("T53").Value = 7.4
$ws.Range("U53").Value = 9
$ws.Range("V53").Value = 8.5
$ws.Range("W53").Value = 16
$ws.Range("X53").Value = 15
$ws.Range("Y53").Value = 27
$ws.Range("Z53").Value = 7.5
$ws.Range("AA53").Value = 6.9
$ws.Range("AB53").Value = 15
$ws.Range("AC53").Value = 70
$ws.Range("AD53").Value = 500
$ws.Range("AE53").Value = 10.75
$ws.Range("AF53").Value = 19
$ws.Range("AG53").Value = 12.5
$ws.Range("AH53").Value = 50
$ws.Range("AI53").Value = 32
$ws.Range("AJ53").Value = 40
$ws.Range("G59").Value = 1.7
$ws.Range("H59").Value = 4.05
$ws.Range("I59").Value = 4.15
$ws.Range("Q59").Value = 3.7
$ws.Range("S59").Value = 2.57
$ws.Range("T59").Value = 11.75
$ws.Range("AE59").Value = 20
$ws.Range("AJ59").Value = 27
